$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for "GET <site>/conseillers/{matricule}/inscriptions" /
# getInscriptionsFromConseiller(String matricule) — it shifts every
# subsequent row up by one and shrinks the table by a row.
$ws.Rows("8:8").Delete()

# Fill in the "Test" (column F) results for the WS verification pass.
$ws.Range("F5").Value = "OK"
$ws.Range("F6").Value = "OK"
$ws.Range("F7").Value = "OK"
$ws.Range("F8").Value = "OK"
$ws.Range("F9").Value = "ne fonctionne pas toujours…"
$ws.Range("F11").Value = "OK"
$ws.Range("F12").Value = "OK"
$ws.Range("F13").Value = "OK"

# Column F widened (bestFit after the new, longer text was added).
$ws.Columns("F:F").ColumnWidth = 23.5

# Selection follows the last edited cell.
$ws.Range("F13").Select()
